# engineer_data.py change: capitalize the hex-byte codes (everything after
# the "0x" prefix) in the "doip" (G) and "uds" (H) columns so lookups are
# case-insensitive-safe / easier to match. Non-hex values (e.g. "N/A") are
# left untouched, and the "0x" prefix itself stays lower-case.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function ConvertHexToUpper($s) {
    if ($null -eq $s) {
        return $s
    }
    $parts = $s.Split(":")
    $newParts = @()
    foreach ($p in $parts) {
        if ($p.Length -gt 2 -and $p.Substring(0,2) -eq "0x") {
            $prefix = $p.Substring(0,2)
            $rest = $p.Substring(2)
            $newParts += $prefix + $rest.ToUpper()
        } else {
            $newParts += $p
        }
    }
    return [string]::Join(":", $newParts)
}

$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    foreach ($col in @("G", "H")) {
        $cell = $ws.Range("$col$row")
        $current = $cell.Text
        $cell.Value = ConvertHexToUpper $current
    }
}
